# "fix typo in included patterns"
#
# Slide 4 ("Used Patterns") has a bulleted list in its content placeholder.
# The first bullet described the scheduler as a "single, globally accessible
# job scheduler" -- fix the wording to "a single, globally accessible job
# queue."

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(4)
$sh = $s.Shapes.Item(2)

if ($sh.Name -ne "Content Placeholder 2") {
    for ($i = 1; $i -le $s.Shapes.Count; $i++) {
        if ($s.Shapes.Item($i).Name -eq "Content Placeholder 2") {
            $sh = $s.Shapes.Item($i)
        }
    }
}

$tr = $sh.TextFrame.TextRange
$para1 = $tr.Paragraphs(1, 1)

$oldPhrase = "accessible job scheduler"
$newPhrase = "accessible job queue."

$offset = $para1.Text.IndexOf($oldPhrase)
if ($offset -ge 0) {
    $target = $tr.Characters($para1.Start + $offset, $oldPhrase.Length)
    $target.Text = $newPhrase
}
